$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new "2022-Q3" row after the header and
#    shift the rest down. Simplest/most-robust way is to just rewrite the
#    whole A2:D9 block with the final values (the index column A is a plain
#    0-based running counter, so it has to be recomputed anyway).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q3", 7,  7.58),
    @(1, "2022-Q2", 7,  9.359999999999999),
    @(2, "2022-Q1", 11, 15.85),
    @(3, "2021-Q4", 26, 23.93),
    @(4, "2021-Q3", 24, 24.2),
    @(5, "2021-Q2", 27, 27.56),
    @(6, "2021-Q1", 39, 25.4),
    @(7, "2020-Q4", 29, 15.41)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# Row 9 is brand new - give its index cell (A9) the same formatting as the
# rest of column A (style index used by A2:A8, bold/centered/bordered).
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)   # xlPasteFormats
$summary.Cells.Item(9, 1).Value = 7

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q3" fund-holdings sheet. Duplicate the existing
#    "2022-Q2" sheet (same column layout) right before it, rename the copy,
#    then overwrite its figures with the 2022-Q3 numbers.
# ---------------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$oldQ2.Copy($oldQ2)

# `.Index` on the handle captured before the copy is stale once the sheet
# collection has been mutated, so re-resolve "2022-Q2" fresh before using
# its position to find the newly-inserted copy sitting right in front of it.
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Item($oldQ2.Index - 1)
$q3.Name = "2022-Q3"

# Helper: write a value as literal text (keeps the same "t=inlineStr / t=s"
# cell type the source data uses for these numeric-looking strings) instead
# of letting it auto-convert to a number, then drop the quote-prefix style
# that Excel tags the cell with so formatting stays identical to its peers.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2 - 005267 嘉实价值精选股票 (code/name unchanged)
Set-TextValue $q3.Range("D2") "46.95"
Set-TextValue $q3.Range("E2") "91.28"
Set-TextValue $q3.Range("F2") "6.09"
Set-TextValue $q3.Range("G2") "2.8593"
$q3.Range("H2").Value = 8

# Row 3 - 010273 嘉实价值长青混合A (code/name unchanged)
Set-TextValue $q3.Range("D3") "40.68"
Set-TextValue $q3.Range("E3") "88.17"
Set-TextValue $q3.Range("F3") "5.31"
Set-TextValue $q3.Range("G3") "2.1601"
$q3.Range("H3").Value = 9

# Row 4 - 001878 嘉实沪港深精选股票 (code/name unchanged)
Set-TextValue $q3.Range("D4") "22.02"
Set-TextValue $q3.Range("E4") "91.33"
Set-TextValue $q3.Range("F4") "4.95"
Set-TextValue $q3.Range("G4") "1.0900"
$q3.Range("H4").Value = 4

# Row 5 - now 160726 嘉实瑞享定期开放灵活配置混合 (was 009138 row)
Set-TextValue $q3.Range("B5") "160726"
Set-TextValue $q3.Range("C5") "嘉实瑞享定期开放灵活配置混合"
Set-TextValue $q3.Range("D5") "12.64"
Set-TextValue $q3.Range("E5") "83.43"
Set-TextValue $q3.Range("F5") "5.13"
Set-TextValue $q3.Range("G5") "0.6484"
$q3.Range("H5").Value = 6

# Row 6 - now 009138 嘉实瑞成两年持有期混合A (was 160726 row)
Set-TextValue $q3.Range("B6") "009138"
Set-TextValue $q3.Range("C6") "嘉实瑞成两年持有期混合A"
Set-TextValue $q3.Range("D6") "10.97"
Set-TextValue $q3.Range("E6") "90.59"
Set-TextValue $q3.Range("F6") "4.23"
Set-TextValue $q3.Range("G6") "0.4640"
$q3.Range("H6").Value = 6

# Row 7 - 010274 嘉实价值长青混合C (code/name unchanged)
Set-TextValue $q3.Range("D7") "4.29"
Set-TextValue $q3.Range("E7") "88.17"
Set-TextValue $q3.Range("F7") "5.31"
Set-TextValue $q3.Range("G7") "0.2278"
$q3.Range("H7").Value = 9

# Row 8 - 009139 嘉实瑞成两年持有期混合C (code/name unchanged)
Set-TextValue $q3.Range("D8") "2.99"
Set-TextValue $q3.Range("E8") "90.59"
Set-TextValue $q3.Range("F8") "4.23"
Set-TextValue $q3.Range("G8") "0.1265"
$q3.Range("H8").Value = 6
